$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the previous layout (B3:I10) entirely before laying out the new tables
$ws.Range("B3:I10").Clear()

# Add the brand-new label text first, in the same order the author typed it,
# so the shared-string table ends up with matching ordering.
$ws.Range("B19").Value = "ingredientsInRes"
$ws.Range("B3").Value = "stored_recipies"
$ws.Range("B10").Value = "stored_ingredients"

# --- Table: stored_recipies ---
$ws.Range("B5").Value = "rec_id"
$ws.Range("C5").Value = "serial"

$ws.Range("B6").Value = "name"
$ws.Range("C6").Value = "varchar"

$ws.Range("B7").Value = "desc"
$ws.Range("C7").Value = "text"

# --- Table: stored_ingredients ---
$ws.Range("B12").Value = "ing_id"
$ws.Range("C12").Value = "serial"

$ws.Range("B13").Value = "name"
$ws.Range("C13").Value = "varchar"

$ws.Range("B14").Value = "unit"
$ws.Range("C14").Value = "varchar"

$ws.Range("B15").Value = "price/unit"
$ws.Range("C15").Value = "double"

$ws.Range("B16").Value = "e/unit"
$ws.Range("C16").Value = "double"

$ws.Range("B17").Value = "p/unit"
$ws.Range("C17").Value = "double"

# --- Table: ingredientsInRes ---
$ws.Range("B21").Value = "rec_id"
$ws.Range("C21").Value = "serial"

$ws.Range("B22").Value = "ing_id"
$ws.Range("C22").Value = "serial"

$ws.Range("B23").Value = "quantity"
$ws.Range("C23").Value = "double"

$ws.Range("B24").Value = "PK(res_id, ing_id"

# Update selection / active cell, matching the author's saved view state
$ws.Range("G9").Select()

# Update window width to reflect the author's resized Excel window
$excel.ActiveWindow.Width = 8268
